$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N8").Value = "FD0810004-ND"
$ws.Range("L8").Value = "FD0810004"
$ws.Range("K8").Value = "Diodes Incorporated"
$ws.Range("H8").Value = "2.97-3.63v"
$ws.Range("I8").Value = "5x3mm"
$ws.Range("O8").Value = "https://www.diodes.com/assets/Datasheets/FD_3-3V.pdf"
$ws.Range("P8").Value = "0Dan_Clock&Timing - Oscillators  - Fixed - Logic_Out:FD0810004"
$ws.Range("Q8").Value = "Osc_Logic_Out:FD0810_DIO"

$ws.Range("E8").Value = "8.192MHz"
$ws.Range("F8").Value = [char]0xB1 + "25PPM"
$ws.Range("G8").Value = "CMOS"
$ws.Range("J8").Value = "XO"
$ws.Range("M8").Value = "Digikey"
$ws.Range("R8").Value = [char]0xA0

$ws.Range("C8").Formula = '=_xlfn.CONCAT(E8,F8," ",H8," ",G8," ",J8," ",I8," ",R8)'
$ws.Range("D8").Formula = '=_xlfn.CONCAT(E8," ",J8)'
$ws.Range("U8").Formula = '=COUNTBLANK(C8:R8)'
$ws.Range("V8").Formula = '=100*COUNTA(C8:R8)/$Z$7'
$ws.Range("W8").Formula = '=IF(V8=100,1,0)'
